$wb = $excel.ActiveWorkbook

# Sheets: 1 = 人物 (People), 2 = 怪物 (Monster), 3 = 特殊 (Special)
$wsPeople  = $wb.Worksheets.Item(1)
$wsMonster = $wb.Worksheets.Item(2)
$wsSpecial = $wb.Worksheets.Item(3)

# The "DropItem" column (P) type row (row 2) changes from int[] to string[]
# on every sheet -- this is the "add drop name column" change: drop items
# are no longer plain int ids but id;name style strings.
$wsPeople.Range("P2").Value = "string[]"
$wsMonster.Range("P2").Value = "string[]"
$wsSpecial.Range("P2").Value = "string[]"

# The People sheet had hard-coded numeric DropItem values left over from
# the old int[] format; clear them out now that the column is string[].
$wsPeople.Range("P4:P42").ClearContents()

# Leave the selection on the just-edited column/row for each sheet, and
# make 特殊 (Special) the active sheet/tab, matching where editing ended.
$wsPeople.Range("P2").Select()
$wsMonster.Range("P2").Select()
$wsSpecial.Range("P3").Select()
$wsSpecial.Activate()
